# "Generate Report for Handoff"
# Refreshes the handoff report: new commit GUID, new content hashes and new
# handoff timestamps, across the Overview, zh-cn and de-de sheets. Each
# hyperlink keeps pointing at its original target URL; only the cell text /
# displayed hyperlink text changes.

$wb = $excel.ActiveWorkbook

$oldGuid = "49466ba0-2874-44b4-bf95-2ef87fd2f651"
$newGuid = "67b8735c-caec-4234-9e09-1e42a1bcf3e5"

$oldHash = "b7eda17806581e498b8ae1b67faacb2c4bd487e6"
$newHash = "edccc569a50ec8d5052c85c6af5379c5bb08cce9"

$newOverviewDate = "2016-45-19 12:45:18"
$newZhDate = "2016-03-19 12:45:16"
$newDeDate = "2016-03-19 12:45:18"

$mdName = "$newGuid.md"
$zhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$deXlfName = "$newGuid.$newHash.de-de.xlf"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/e0c5620183ff90da78cecf923d4eb10e0ad10da4/e2e/$oldGuid.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bbdfa34bd600a6956ebd6e9d312635b6eb1e9b18/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bdcd80cce783331d887cb08fb96ec4b89e8b4998/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf"

# ---- Overview sheet: A2 hyperlink display, D2 handoff date ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, "", "", $mdName)
$ws.Range("D2").Value = $newOverviewDate

# ---- zh-cn sheet: A2/B2/D2 hyperlink displays, E2 handoff datetime ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdAddress, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $zhXlfAddress, "", "", $zhXlfName)
$ws.Range("E2").Value = $newZhDate

# ---- de-de sheet: A2/B2/D2 hyperlink displays, E2 handoff datetime ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, "", "", $mdName)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdAddress, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $deXlfAddress, "", "", $deXlfName)
$ws.Range("E2").Value = $newDeDate

Write-Host "Handoff report regenerated."
